$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "RGossF"

# 2. Fix tiny floating point differences in existing cells
$ws.Range("H13").Value = 0.9990324650035259
$ws.Range("E15").Value = 0.9778202286087941
$ws.Range("H15").Value = 0.9697320750700568
$ws.Range("M15").Value = 1.001700669685752

# 3. Add new row 16 of data
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9834174039073463
$ws.Range("D16").Value = 0.9924703951312377
$ws.Range("E16").Value = 1.027488790855159
$ws.Range("F16").Value = 0.9834174039073463
$ws.Range("G16").Value = 0.9211094504501502
$ws.Range("H16").Value = 1.201031381252061
$ws.Range("I16").Value = 1.007901537378366
$ws.Range("J16").Value = 0.9924703951312377
$ws.Range("K16").Value = 1.009979592993199
$ws.Range("L16").Value = 0.9966984984502724
$ws.Range("M16").Value = 1.022236493162387
